$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, shifting existing rows 97:210 down to 98:211
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with fresh data (columns that stay
# constant across the sheet are copied from the template row immediately
# below, which now lives at row 98)
$ws.Range("A97").Value = 5
$ws.Range("B97").Value = "Macroferia Regional de Talca"
$ws.Range("C97").Value = "Maule"
$ws.Range("D97").Value = 45174
$ws.Range("D97").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E97").Value = 7
$ws.Range("F97").Value = 100112001
$ws.Range("G97").Value = "Berenjena"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = 10000
$ws.Range("N97").Value = "$/caja 50 unidades"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 200
$ws.Range("Q97").Value = 50
$ws.Range("R97").Value = "Hortaliza"
